$d = $word.ActiveDocument

# Locate the paragraph that currently ends the document ("fuck you paul")
# and build an insertion point immediately after its paragraph mark, i.e.
# at the very end of the document's content.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertAt = $lastPara.Range.End
$ip = $d.Range($insertAt, $insertAt)

# Insert two new paragraphs after it:
#   1) an empty paragraph (Times New Roman end-of-paragraph mark formatting)
#   2) a paragraph containing "i love you man" (also Times New Roman)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:pPr>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
                '</w:rPr>' +
              '</w:pPr>' +
            '</w:p>' +
            '<w:p>' +
              '<w:pPr>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
                '</w:rPr>' +
              '</w:pPr>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
                '</w:rPr>' +
                '<w:t>i love you man</w:t>' +
              '</w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$ip.InsertXML($xml)
